$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in column B
$ws.Range("B2").Value = 221
$ws.Range("B3").Value = 206

# Add new rows 4 and 5
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 105

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 98

# Copy the style (bold, border, centered) from A3 to the new A4/A5 cells
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
